$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry")

# Update delta pH formulas for the 3.4 kg/t dose rows (D4 and D8)
$ws.Range("D4").Formula = "=7.9-0.8187"
$ws.Range("D8").Formula = "=7.9-0.8187"

# Update the active cell selection on the Slurry sheet
$ws.Activate()
$ws.Range("E8").Select()
